$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").Value = "sex"

$ws.Range("J2").Value = "Male"
$ws.Range("J3").Value = "Female"
$ws.Range("J4").Value = "Female"
$ws.Range("J5").Value = "Male"
$ws.Range("J6").Value = "Female"
$ws.Range("J7").Value = "Male"
$ws.Range("J8").Value = "Male"
$ws.Range("J9").Value = "Female"
$ws.Range("J10").Value = "Female"
$ws.Range("J11").Value = "Female"
$ws.Range("J12").Value = "Female"
$ws.Range("J13").Value = "Male"
$ws.Range("J14").Value = "Female"
$ws.Range("J15").Value = "Female"
$ws.Range("J16").Value = "Male"
$ws.Range("J17").Value = "Female"
$ws.Range("J18").Value = "Female"
$ws.Range("J19").Value = "Female"
$ws.Range("J20").Value = "Female"
$ws.Range("J21").Value = "Male"
$ws.Range("J22").Value = "Male"
$ws.Range("J23").Value = "Male"
$ws.Range("J24").Value = "Male"
$ws.Range("J25").Value = "Male"
$ws.Range("J26").Value = "Male"
$ws.Range("J27").Value = "Male"
$ws.Range("J28").Value = "Male"
$ws.Range("J29").Value = "Male"
$ws.Range("J30").Value = "Male"
$ws.Range("J31").Value = "Female"
$ws.Range("J32").Value = "Male"
$ws.Range("J33").Value = "Male"
$ws.Range("J34").Value = "Female"
$ws.Range("J35").Value = "Male"
$ws.Range("J36").Value = "Female"
$ws.Range("J37").Value = "Male"
$ws.Range("J38").Value = "Female"
$ws.Range("J39").Value = "Male"
$ws.Range("J40").Value = "Female"
$ws.Range("J41").Value = "Male"
$ws.Range("J42").Value = "Female"
$ws.Range("J43").Value = "Male"
$ws.Range("J44").Value = "Male"
$ws.Range("J45").Value = "Female"
$ws.Range("J46").Value = "Male"
$ws.Range("J47").Value = "Male"
$ws.Range("J48").Value = "Male"
$ws.Range("J49").Value = "Male"
$ws.Range("J50").Value = "Female"
$ws.Range("J51").Value = "Male"
$ws.Range("J52").Value = "Female"
$ws.Range("J53").Value = "Female"
$ws.Range("J54").Value = "Female"
$ws.Range("J55").Value = "Male"
$ws.Range("J56").Value = "Male"
$ws.Range("J57").Value = "Female"
$ws.Range("J58").Value = "Female"
$ws.Range("J59").Value = "Female"
$ws.Range("J60").Value = "Male"
$ws.Range("J61").Value = "Male"
$ws.Range("J62").Value = "Female"
$ws.Range("J63").Value = "Female"
$ws.Range("J64").Value = "Female"
$ws.Range("J65").Value = "Male"
$ws.Range("J66").Value = "Male"
$ws.Range("J67").Value = "Female"
$ws.Range("J68").Value = "Female"
$ws.Range("J69").Value = "Female"
$ws.Range("J70").Value = "Male"
$ws.Range("J71").Value = "Male"
$ws.Range("J72").Value = "Female"
$ws.Range("J73").Value = "Female"
$ws.Range("J74").Value = "Male"
$ws.Range("J75").Value = "Male"
$ws.Range("J76").Value = "Male"
$ws.Range("J77").Value = "Male"
$ws.Range("J78").Value = "Male"
$ws.Range("J79").Value = "Male"
$ws.Range("J80").Value = "Male"
$ws.Range("J81").Value = "Male"
$ws.Range("J82").Value = "Female"
$ws.Range("J83").Value = "Female"
$ws.Range("J84").Value = "Female"
$ws.Range("J85").Value = "Male"
$ws.Range("J86").Value = "Male"
$ws.Range("J87").Value = "Female"
$ws.Range("J88").Value = "Male"
$ws.Range("J89").Value = "Female"
$ws.Range("J90").Value = "Male"
$ws.Range("J91").Value = "Male"
$ws.Range("J92").Value = "Female"
$ws.Range("J93").Value = "Female"
$ws.Range("J94").Value = "Female"
$ws.Range("J95").Value = "Female"
$ws.Range("J96").Value = "Male"
$ws.Range("J97").Value = "Female"
$ws.Range("J98").Value = "Male"
$ws.Range("J99").Value = "Female"
$ws.Range("J100").Value = "Male"
$ws.Range("J101").Value = "Male"
$ws.Range("J102").Value = "Male"
$ws.Range("J103").Value = "Male"
$ws.Range("J104").Value = "Male"
$ws.Range("J105").Value = "Male"
$ws.Range("J106").Value = "Male"
$ws.Range("J107").Value = "Male"
$ws.Range("J108").Value = "Male"
$ws.Range("J109").Value = "Male"
$ws.Range("J110").Value = "Male"
$ws.Range("J111").Value = "Male"
$ws.Range("J112").Value = "Male"
$ws.Range("J113").Value = "Male"
$ws.Range("J114").Value = "Male"
$ws.Range("J115").Value = "Female"
$ws.Range("J116").Value = "Female"
$ws.Range("J117").Value = "Male"
$ws.Range("J118").Value = "Female"
$ws.Range("J119").Value = "Male"
$ws.Range("J120").Value = "Female"

$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("M23").Select()
